# Change the policy of output YAML file.
# Each sheet included in excel file will be written by each sheet.
#
# This adds a new "Note" sheet (free-form notes / title / update-date /
# comment) in front of the existing "device" sheet, and promotes the
# "device" sheet's former generic placeholder labels ("Title:",
# "Update Date:", "Comment:") to the actual, specific values that used
# to live only implicitly. The "Note" sheet keeps the generic labels as
# its own header rows, plus a free-format note area and a trailing
# disclaimer comment.

$wb = $excel.ActiveWorkbook
$device = $wb.ActiveSheet

# --- Insert a new sheet "Note" in front of the "device" sheet --------
$note = $wb.Worksheets.Add($device)
$note.Name = "Note"

# --- Populate the "Note" sheet ---------------------------------------
# Generic header labels (reuse the same text/styles the "device" sheet
# used to show before this edit).
$note.Range("A1").Value = "Title:"
$note.Range("A1").Font.Bold = $true

$note.Range("A2").Value = "Update Date:"
$note.Range("A2").Font.Bold = $true

$note.Range("A3").Value = "Comment:"
$note.Range("A3").Font.Bold = $true
$note.Range("A3").VerticalAlignment = -4160

# Free-format blank area (rows 4-25, columns A-I) - just touch the
# alignment so the cells materialize with a (blank) style, matching a
# "this area is free-form" visual block.
$note.Range("A4:I25").WrapText = $false

# Trailing disclaimer comment (row 26 intentionally left empty).
$note.Range("A27").Value = "* This sheet is free format and will not be transformed to YAML file."

# --- Update the "device" sheet with the real title/date/comment ------
$device.Range("A1").Value = "Title: [TEST 3-1-1] Single chassis topology"
$device.Range("A2").Value = "Update Date: 1984/09/01 by yuji"
$device.Range("A3").Value = "Comment: Single chassis topology at B1 DC"

# Selection on "device" now spans the header rows (1-3, all columns).
$device.Range("A1:XFD3").Select()

# "device" stays the active tab (second tab, index 1).
$wb.Worksheets.Item(2).Activate()

# Printable page setup for the "device" sheet.
$device.PageSetup.PaperSize = 9
$device.PageSetup.Orientation = 1
